$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.934.96'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '3.508.56'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'" + '595.91'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = "'" + '169.80'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = "'" + '0.591'
$ws.Range('E8').Value = '  +3.50%  '
$ws.Range('D9').Value = "'" + '0.134'
$ws.Range('E9').Value = '  +7.66%  '
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '4.116.69'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('E14').Value = '  +1.94%  '
$ws.Range('D15').Value = "'" + '0.0000182'
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('D16').Value = '66.922.10'
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = '3.487.33'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').Value = "'" + '395.98'
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('D21').Value = "'" + '8.01'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = "'" + '73.63'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('E24').Value = '  +2.20%  '
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('D26').Value = "'" + '10.23'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = "'" + '0.998'
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').Value = "'" + '6.31'
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('D30').Value = "'" + '1.46'
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = "'" + '24.00'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').Value = "'" + '7.40'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').Value = "'" + '1.61'
$ws.Range('E34').Value = '  +3.89%  '
$ws.Range('D35').Value = "'" + '163.97'
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').Value = "'" + '1.92'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = "'" + '4.75'
$ws.Range('E38').Value = '  +4.19%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = "'" + '6.88'
$ws.Range('E39').Value = '  +2.78%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('D42').Value = "'" + '2.64'
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.830.19'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = "'" + '27.06'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = "'" + '42.88'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').Value = "'" + '0.0313'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').Value = "'" + '342.38'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('D48').Value = "'" + '1.11'
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('D49').Value = "'" + '33.85'
$ws.Range('E49').Value = '  +2.91%  '
$ws.Range('E50').Value = '  +0.83%  '
$ws.Range('E51').Value = '  +1.32%  '
